$d = $word.ActiveDocument

function Assert-ParaText($index, $expectedPrefix) {
    $actual = $d.Paragraphs.Item($index).Range.Text
    if (-not $actual.StartsWith($expectedPrefix)) {
        throw "Paragraph $index mismatch: expected prefix [$expectedPrefix] got [$actual]"
    }
}

function Format-HeaderParagraph($para) {
    # Bold the paragraph's visible text (but not the trailing paragraph
    # mark) and set before/after spacing to 6pt (120 dxa / 120 twips).
    # This reproduces the target:
    #   <w:pPr><w:spacing w:before="120" w:after="120"/></w:pPr>
    #   <w:r><w:rPr><w:b/></w:rPr>...
    # without also bolding the paragraph mark itself (which would add an
    # unwanted <w:rPr><w:b/></w:rPr> inside <w:pPr>).
    $r = $para.Range
    $textOnly = $d.Range($r.Start, $r.End - 1)
    $textOnly.Bold = 1
    $para.Format.SpaceBefore = 6
    $para.Format.SpaceAfter = 6
}

# ---------------------------------------------------------------------------
# Sanity-check the paragraph layout we expect to find in the original doc
# before mutating anything.
# ---------------------------------------------------------------------------
Assert-ParaText 2  "[EM] Consolidation"
Assert-ParaText 3  "2. Customer Research"
Assert-ParaText 4  "3. Our Solution"
Assert-ParaText 7  "Dependencies:"
Assert-ParaText 8  "4. Product Metrics"
Assert-ParaText 13 "Appendix: LinksHYPERLINK"
Assert-ParaText 14 "Appendix: Quick prototype"
Assert-ParaText 15 "Figure: PDF page 1"

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so the paragraph indices used
# below stay valid as new paragraphs get inserted above them. Within each
# block, new blank/plain paragraphs are inserted first (by calling
# InsertParagraphAfter/Before on a paragraph that is still plain "Normal"
# formatting) BEFORE the section-header paragraph itself is bolded/spaced -
# this way the new blank/plain paragraphs never inherit the bold, spacing,
# or heading-style formatting that gets applied to the header.
# ---------------------------------------------------------------------------

# === "Appendix: Quick prototype" (Heading2) =================================
$pProto = $d.Paragraphs.Item(14)
$pFigure = $d.Paragraphs.Item(15)             # "Figure: PDF page 1" - plain Normal
$pFigure.Range.InsertParagraphBefore()        # blank line after the header; inherits Normal from pFigure, not Heading2 from pProto
Format-HeaderParagraph $pProto

# === blank line after "Appendix: Links...HYPERLINK..." ======================
$pLinks = $d.Paragraphs.Item(13)
$pLinks.Range.InsertParagraphAfter()

# === "4. Product Metrics..." =================================================
$pMetrics = $d.Paragraphs.Item(8)
$pMetrics.Range.InsertParagraphAfter()
Format-HeaderParagraph $pMetrics

# === blank line after "Dependencies: ..." ====================================
$pDeps = $d.Paragraphs.Item(7)
$pDeps.Range.InsertParagraphAfter()

# === "3. Our Solution..." =====================================================
$pSolution = $d.Paragraphs.Item(4)
$pSolution.Range.InsertParagraphAfter()
Format-HeaderParagraph $pSolution

# === "2. Customer Research..." + two new body paragraphs =====================
$pResearch = $d.Paragraphs.Item(3)
$pResearch.Range.InsertParagraphAfter()       # blank line right after the header (index 4)

$pBlank1 = $d.Paragraphs.Item(4)
$pBlank1.Range.InsertParagraphAfter()         # room for first new paragraph (index 5)
$pNew1 = $d.Paragraphs.Item(5)
$pNew1.Range.Text = "This capability was requested as feedback from an enterprise-level accounting firm, reflecting needs observed in large multi-entity audit workflows."

$pNew1.Range.InsertParagraphAfter()           # room for second new paragraph (index 6)
$pNew2 = $d.Paragraphs.Item(6)
$pNew2.Range.Text = "We are also building this to achieve competitive parity with Wolters Kluwer ProSystem fx Engagement, which offers similar functionality."

$pNew2.Range.InsertParagraphAfter()           # trailing blank line (index 7)

Format-HeaderParagraph $pResearch             # bold + spacing applied last so it never leaks onto 4-7

# === blank line after "[EM] Consolidation..." ================================
$pIntro = $d.Paragraphs.Item(2)
$pIntro.Range.InsertParagraphAfter()

Write-Output "done: paragraphs=$($d.Paragraphs.Count)"
